$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text cells: values that look numeric must stay as text (matches source formatting)
$ws.Range("D2").Value = '27.563.01'
$ws.Range("E2").Value = '  -0.88%  '
$ws.Range("D3").Value = '1.662.29'
$ws.Range("E3").Value = '  -3.52%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.74'
$ws.Range("E5").Value = '  -1.02%  '
$ws.Range("E6").Value = '  -1.95%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.37'
$ws.Range("E8").Value = '  -2.12%  '
$ws.Range("E9").Value = '  -2.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0621'
$ws.Range("E10").Value = '  -1.38%  '
$ws.Range("E11").Value = '  -2.62%  '
$ws.Range("D12").Value = '1.896.93'
$ws.Range("E12").Value = '  -3.59%  '
$ws.Range("D13").Value = '1.657.45'
$ws.Range("E13").Value = '  -3.66%  '
$ws.Range("E14").Value = '  -2.69%  '
$ws.Range("E15").Value = '  -3.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.90'
$ws.Range("E16").Value = '  -2.99%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '245.98'
$ws.Range("E17").Value = '  +2.08%  '
$ws.Range("D18").Value = '27.531.73'
$ws.Range("E18").Value = '  -1.10%  '
$ws.Range("D19").Value = '0.0₃0732'
$ws.Range("E19").Value = '  -2.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.43'
$ws.Range("E20").Value = '  -7.66%  '
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("E22").Value = '  -3.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.33'
$ws.Range("E23").Value = '  -3.60%  '
$ws.Range("E24").Value = '  -4.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.09'
$ws.Range("E25").Value = '  -1.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.16'
$ws.Range("E26").Value = '  -4.92%  '
$ws.Range("E27").Value = '  -2.44%  '
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.111'
$ws.Range("E29").Value = '  -2.22%  '
$ws.Range("E30").Value = '  +4.39%  '
$ws.Range("E31").Value = '  -1.23%  '
$ws.Range("E32").Value = '  -3.08%  '
$ws.Range("D33").Value = '1.447.58'
$ws.Range("E33").Value = '  -1.55%  '
$ws.Range("E34").Value = '  -5.20%  '
$ws.Range("E35").Value = '  -7.33%  '
$ws.Range("E36").Value = '  -0.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.931'
$ws.Range("E37").Value = '  -3.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.579'
$ws.Range("E38").Value = '  -5.69%  '
$ws.Range("E39").Value = '  -2.71%  '
$ws.Range("E40").Value = '  -3.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.09'
$ws.Range("E41").Value = '  -3.56%  '
$ws.Range("E42").Value = '  +0.21%  '
$ws.Range("E43").Value = '  -7.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.790'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = '1.805.03'
$ws.Range("E45").Value = '  -3.49%  '
$ws.Range("E46").Value = '  -3.49%  '
$ws.Range("E47").Value = '  +0.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.75'
$ws.Range("E48").Value = '  -3.24%  '
$ws.Range("D49").Value = '0.0₆0108'
$ws.Range("E49").Value = '  -1.13%  '
$ws.Range("E50").Value = '  -4.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.85'
$ws.Range("E51").Value = '  -5.20%  '
